# This script updates the "想去人数" (want-to-go headcount) figures in
# column F across the 展览, 演出, and 全部类型 worksheets, matching the
# newly regenerated data output (commit "Update gh-pages to output
# generated at 456a3b4").

$wb = $excel.ActiveWorkbook

function Set-FValues {
    param(
        $SheetName,
        $RowValues
    )

    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($row in $RowValues.Keys) {
        $ws.Range("F$row").Value = $RowValues[$row]
    }
}

# Sheet: 展览
Set-FValues "展览" @{
    3  = 1024
    4  = 806
    5  = 884
    6  = 462
    7  = 710
    9  = 1312
    12 = 560
    15 = 1145
    16 = 139
    18 = 424
    19 = 379
    21 = 598
    22 = 157
    23 = 663
    25 = 1073
}

# Sheet: 演出
Set-FValues "演出" @{
    2  = 347
    7  = 253
    11 = 115
}

# Sheet: 全部类型
Set-FValues "全部类型" @{
    3  = 347
    5  = 1024
    6  = 806
    7  = 884
    8  = 462
    9  = 462
    10 = 710
    12 = 1312
    17 = 560
    21 = 1146
    23 = 139
    25 = 424
    26 = 379
    28 = 253
    30 = 598
    33 = 115
    34 = 115
    35 = 157
    36 = 663
    38 = 1073
}
